# EWD-22592 - Sign-Up in 2 steps buttons
# Fill in missing NL/DE translations in the localization resources sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# savedAt
$ws.Range("D94").Value = "Opgeslagen om"
$ws.Range("E94").Value = "Gespeichert am"

# Server error notification block
$ws.Range("D211").Value = "De leerervaring is niet gevonden. Vernieuw uw gegevens. "
$ws.Range("E211").Value = "Die Lernerfahrung wurde nicht gefunden. Bitte aktualisieren Sie Ihre Daten."

$ws.Range("D212").Value = "Het leerdoel is niet gevonden. Vernieuw uw gegevens."
$ws.Range("E212").Value = "Das Lernziel wurde nicht gefunden. Bitte aktualisieren Sie Ihre Daten."

$ws.Range("D213").Value = "De antwoordoptie is niet gevonden. Vernieuw uw gegevens."
$ws.Range("E213").Value = "Die Antwortoption wurde nicht gefunden. Bitte aktualisieren Sie Ihre Daten."

$ws.Range("D214").Value = "De vraag is niet gevonden. Vernieuw uw gegevens."
$ws.Range("E214").Value = "Die Frage wurde nicht gefunden. Bitte aktualisieren Sie Ihre Daten."

$ws.Range("D215").Value = "Het leerobject is niet gevonden. Vernieuw uw gegevens."
$ws.Range("E215").Value = "Der Lerngegenstand wurde nicht gefunden. Bitte aktualisieren Sie Ihre Daten."

$ws.Range("D216").Value = "De helphint is niet gevonden. Vernieuw uw gegevens."
$ws.Range("E216").Value = "Der Hilfehinweis wurde nicht gefunden. Bitte aktualisieren Sie Ihre Daten."

# View captions block
$ws.Range("D218").Value = "Leerervaring-editor"

$ws.Range("D219").Value = "Leerdoel-editor"
$ws.Range("E219").Value = "Lernziel-Editor"

$ws.Range("D220").Value = "Vraag-editor"
$ws.Range("E220").Value = "Frage-Editor"

$ws.Range("D221").Value = "Leerervaring maken"
$ws.Range("E221").Value = "Lernerfahrung erstellen"

$ws.Range("D223").Value = "Vraag maken"
$ws.Range("E223").Value = "Frage erstellen"

# Sign-UP optional fields block
$ws.Range("D225").Value = "Volledige naam"
$ws.Range("E225").Value = "Vollständiger Name"

$ws.Range("D226").Value = "Telefoonnummer"
$ws.Range("E226").Value = "Telefonnummer"

$ws.Range("D227").Value = "Naam van organisatie"
$ws.Range("E227").Value = "Name der Organisiation"

$ws.Range("D228").Value = "Land"

$ws.Range("D229").Value = "Hoeveel medewerkers in uw bedrijf schrijven cursusmateriaal?"
$ws.Range("E229").Value = "Wie viele Personen in Ihrem Unternehmen erstellen Kurse?"

$ws.Range("D230").Value = "Hoe snel hebt u een schrijftool nodig?"
$ws.Range("E230").Value = "Wie bald benötigen Sie ein Autorentool?"

$ws.Range("D231").Value = "Hoe schrijft u nu cursusmateriaal?"
$ws.Range("E231").Value = "Wie erstellen Sie heute Kurse?"

$ws.Range("D232").Value = "Kies een antwoord"
$ws.Range("E232").Value = "Antwort auswählen"

$ws.Range("D233").Value = "Nu"
$ws.Range("E233").Value = "Jetzt"

$ws.Range("D234").Value = "1-3 maanden"
$ws.Range("E234").Value = "In 1 - 3 Monaten"

$ws.Range("D235").Value = ">3 maanden"
$ws.Range("E235").Value = "In mehr als 3 Monaten"

$ws.Range("D236").Value = "Ik oriënteer me alleen"
$ws.Range("E236").Value = "Ich informiere mich nur."

$ws.Range("D238").Value = "Een andere schrijftool"
$ws.Range("E238").Value = "Mit einem anderen Autorentool"

$ws.Range("D239").Value = "Geen"
$ws.Range("E239").Value = "Keines"

$ws.Range("D240").Value = "Anders"
$ws.Range("E240").Value = "Sonstiges"

$ws.Range("D241").Value = "Overslaan"
$ws.Range("E241").Value = "Überspringen"

$ws.Range("D242").Value = "Verplichte informatie"
$ws.Range("E242").Value = "Erforderliche Information"

$ws.Range("D243").Value = "Optionele informatie"
$ws.Range("E243").Value = "Optionale Information"

$ws.Range("D244").Value = "Kies uw land"
$ws.Range("E244").Value = "Wählen Sie Ihr Land aus."
